$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (D) and Volume(1h) (E) columns with latest crypto data
$data = @(
    @{Row=2; D="51.581.72"; E="  +1.20%  "},
    @{Row=3; D="2.989.81"; E="  +1.51%  "},
    @{Row=4; D="1.00"; E="  -0.04%  "},
    @{Row=5; D="382.24"; E="  +1.77%  "},
    @{Row=6; D="104.60"; E="  +3.39%  "},
    @{Row=7; E="  +1.24%  "},
    @{Row=8; E="  +0.00%  "},
    @{Row=9; D="0.599"; E="  +1.72%  "},
    @{Row=10; D="36.72"; E="  +1.08%  "},
    @{Row=11; E="  -0.66%  "},
    @{Row=12; D="0.0860"; E="  +1.03%  "},
    @{Row=13; D="3.454.91"; E="  +1.59%  "},
    @{Row=14; D="7.88"; E="  +3.43%  "},
    @{Row=15; E="  +2.22%  "},
    @{Row=16; D="2.988.55"; E="  +1.73%  "},
    @{Row=17; D="11.19"; E="  -0.91%  "},
    @{Row=18; D="0.996"; E="  -0.12%  "},
    @{Row=19; D="51.604.23"; E="  +1.29%  "},
    @{Row=20; D="3.09"; E="  +0.38%  "},
    @{Row=21; D="12.60"; E="  +1.04%  "},
    @{Row=22; E="  +0.79%  "},
    @{Row=23; D="70.49"},
    @{Row=24; D="267.79"; E="  +0.52%  "},
    @{Row=25; D="3.24"; E="  +1.69%  "},
    @{Row=26; D="7.91"; E="  -2.75%  "},
    @{Row=27; D="7.29"; E="  -1.69%  "},
    @{Row=28; D="0.169"; E="  +4.06%  "},
    @{Row=29; E="  +0.11%  "},
    @{Row=30; D="26.15"; E="  +1.78%  "},
    @{Row=31; E="  -0.46%  "},
    @{Row=32; E="  +4.42%  "},
    @{Row=33; D="34.70"; E="  +4.03%  "},
    @{Row=34; D="51.44"; E="  +0.96%  "},
    @{Row=35; E="  +0.50%  "},
    @{Row=36; D="0.0448"; E="  +1.35%  "},
    @{Row=37; E="  +0.01%  "},
    @{Row=38; E="  +5.04%  "},
    @{Row=39; D="16.94"; E="  +2.63%  "},
    @{Row=40; D="2.59"; E="  +4.85%  "},
    @{Row=41; E="  +1.02%  "},
    @{Row=42; E="  +2.18%  "},
    @{Row=43; D="3.87"; E="  +13.24%  "},
    @{Row=44; D="126.75"; E="  +5.46%  "},
    @{Row=45; D="21.49"; E="  +1.47%  "},
    @{Row=46; E="  -0.07%  "},
    @{Row=47; D="0.272"; E="  +0.31%  "},
    @{Row=48; E="  +0.88%  "},
    @{Row=49; D="2.038.47"; E="  +2.28%  "},
    @{Row=50; D="3.281.83"; E="  +1.38%  "},
    @{Row=51; E="  +0.38%  "}
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.ContainsKey("D")) {
        # Force text format so numeric-looking strings (e.g. "1.00") are not
        # auto-converted to numbers, preserving original text representation
        $ws.Cells.Item($r, 4).NumberFormat = "@"
        $ws.Cells.Item($r, 4).Value = $item.D
    }
    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}
